$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new record row (row 12) - "YOLLY DEL VALLE DIAZ DIAZ" / Cedula / 17571916 / OFDA
$ws.Range("A12").Value = "YOLLY"
$ws.Range("B12").Value = "DEL VALLE"
$ws.Range("C12").Value = "DIAZ"
$ws.Range("D12").Value = "DIAZ"
$ws.Range("E12").Value = "Cedula"
$ws.Range("F12").Value = 17571916
$ws.Range("G12").Value = "OFDA"

# Move selection to G8, matching the author's last active cell before saving
$ws.Range("G8").Select()
